$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# The commit inserts the honorific "Pan"/"Pana" immediately before each
# of the four standalone mentions of "Karol"/"Karola" in the Polish
# recommendation-letter body (the English copy further down is left
# untouched, matching the diff).
#
# This engine's Range mutation ops (InsertBefore / Text= / Find-replace)
# normalize ("flatten") every run in the touched paragraph down to a
# single run sharing the touched run's formatting before re-inserting -
# that's fine for getting the *text* right, but it also silently
# re-fuses unrelated, pre-existing run boundaries elsewhere in the same
# paragraph (e.g. "od podstaw " / "pracy" used to be their own runs for
# no textual reason - just an artifact of previous edits/rsids).
#
# Formatting-only ops (Font.Bold = ...) on a sub-range, by contrast, are
# surgical: they only ever split the run(s) that sub-range overlaps and
# leave every other run boundary in the paragraph alone. So the recipe
# is: (1) do the plain-text insertions first with InsertBefore, which
# gives correct text but over-merged runs, then (2) make a second pass
# that re-imposes every run boundary the target XML needs (both the
# brand new "Pan"/"Pana" splits AND the pre-existing ones that must
# survive) by toggling Bold on/off across each exact sub-range - a
# bold-then-unbold round trip leaves formatting unchanged but forces a
# run split at both ends of the range.
# ---------------------------------------------------------------------

function Split-Range($startPos, $endPos) {
    if ($endPos -le $startPos) { return }
    $rng = $d.Range($startPos, $endPos)
    $rng.Font.Bold = $true
    $rng.Font.Bold = $false
}

function Insert-TextBefore($searchText, $insertText, $searchStart) {
    $scope = $d.Range($searchStart, $d.Content.End)
    $found = $scope.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NOT FOUND: $searchText"
        return -1
    }
    $pos = $scope.Start
    $ins = $d.Range($pos, $pos)
    $ins.InsertBefore($insertText)
    return $pos
}

function FindStart($searchText, $searchStart) {
    $scope = $d.Range($searchStart, $d.Content.End)
    $found = $scope.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NOT FOUND: $searchText"
        return -1
    }
    return $scope.Start
}

function FindEnd($searchText, $searchStart) {
    $scope = $d.Range($searchStart, $d.Content.End)
    $found = $scope.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NOT FOUND: $searchText"
        return -1
    }
    return $scope.End
}

# =======================================================================
# Phase 1 - fix the text content (four "Pan "/"Pana " insertions)
# =======================================================================

$p1 = Insert-TextBefore "Karol zawsze przychodził" "Pan " 0
$p2 = Insert-TextBefore "Karol pracował samodzielnie" "Pan " 0
$p3 = Insert-TextBefore "Karola było samodzielne" "Pana " 0
$p4 = Insert-TextBefore "Karol aktywnie uczestniczył" "Pan " 0

Write-Output "Phase 1 (text) complete: p1=$p1 p2=$p2 p3=$p3 p4=$p4"

# =======================================================================
# Phase 2 - re-impose the exact run boundaries the target XML expects
# =======================================================================

# --- Paragraph: "...9 sierpnia 2024. [Pan ][Karol zawsze przychodził...]"
$anchorStart = FindEnd "praktyki" 0
$panStart = FindStart "Pan " $anchorStart
$panEnd = $panStart + 4
Split-Range $anchorStart $panStart
Split-Range $panStart $panEnd
$karolEnd = FindEnd "Karol zawsze przychodził do pracy punktualnie i w odpowiednim " $panEnd
Split-Range $panEnd $karolEnd

# --- Paragraph: "Podczas praktyk[ ][Pan ][Karol pracował samodzielnie...]"
$anchorStart = FindStart "praktyk" $karolEnd
$praktykEnd = $anchorStart + 7   # "praktyk" (7 chars) stays its own pre-existing run
$panStart = FindStart "Pan " $praktykEnd
$panEnd = $panStart + 4
Split-Range $praktykEnd $panStart
Split-Range $panStart $panEnd
$restEnd = FindEnd "Karol pracował samodzielnie, wymagając minimalnego nadzoru. Szybko zapoznał się z naszym zestawem technologii, co pozwoliło mu na efektywne wykonywanie zadań." $panEnd
Split-Range $panEnd $restEnd

# --- Paragraph: "Głównym osiągnięciem [Pan][a][ ][Karola było samodzielne stworzenie ]od podstaw ...planowanie [pracy][. ][Pan ][Karol aktywnie...]"
$osStart = FindStart "Głównym osiągnięciem " $restEnd
$osEnd = FindEnd "Głównym osiągnięciem " $restEnd
$panStart = $osEnd
$panMid = $panStart + 3     # end of "Pan"
$aEnd = $panMid + 1         # end of "a"
$spEnd = $aEnd + 1          # end of " "
Split-Range $osStart $osEnd
Split-Range $panStart $panMid
Split-Range $panMid $aEnd
Split-Range $aEnd $spEnd

$stworzenieEnd = FindEnd "Karola było samodzielne stworzenie " $spEnd
Split-Range $spEnd $stworzenieEnd

$odPodstawStart = FindStart "od podstaw " $stworzenieEnd
$odPodstawEnd = FindEnd "od podstaw " $stworzenieEnd
Split-Range $stworzenieEnd $odPodstawStart
Split-Range $odPodstawStart $odPodstawEnd

$planowanieEnd = FindEnd "nowej wersji aplikacji do koordynacji wysyłek pomiędzy bramą, magazynem a logistyką. Nowa wersja, używana przez (jeszcze się dopytam ile) osób, jest nie tylko szybsza, ale również zawiera dodatkowe funkcjonalności, które pozwoliły na bardziej elastyczne planowanie " $odPodstawEnd
Split-Range $odPodstawEnd $planowanieEnd

$pracyStart = FindStart "pracy" $planowanieEnd
$pracyEnd = $pracyStart + 5
Split-Range $planowanieEnd $pracyStart
Split-Range $pracyStart $pracyEnd

$kropkaEnd = $pracyEnd + 2   # ". " (2 chars)
Split-Range $pracyEnd $kropkaEnd

$panStart2 = FindStart "Pan " $kropkaEnd
$panEnd2 = $panStart2 + 4
Split-Range $kropkaEnd $panStart2
Split-Range $panStart2 $panEnd2

$finalEnd = FindEnd "Karol aktywnie uczestniczył w spotkaniach konsultacyjnych, a także samodzielnie kontaktował się z użytkownikami, aby stworzyć jak najlepszą aplikację odpowiadającą ich potrzebom." $panEnd2
Split-Range $panEnd2 $finalEnd

Write-Output "Phase 2 (run boundaries) complete"
